$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme
$cs.Colors(1).RGB = 5649426
